$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text values
$ws.Range("A1").Value = "No"
$ws.Range("B1").Value = "Product-Code"
$ws.Range("C1").Value = "Product-Name"

# Update view: zoom and selection
$excel.ActiveWindow.Zoom = 161
$ws.Range("C3").Select()
